# Temperate Forest Traits.xlsx — fill in newly measured Width/Cross-Sectional-Area
# style paired values (columns E/F) for the individuals in rows 240-271 and
# 306-318 that previously had no measurements recorded. A couple of
# individuals (rows 240 and 256) were not measurable, so they are recorded
# as "X" instead of a number, matching the rest of the workbook's convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(240, "X", "X"),
    @(241, 87, 112),
    @(242, 94, 65),
    @(243, 142, 107),
    @(244, 83, 53),
    @(245, 135, 60),
    @(246, 131, 73),
    @(247, 39, 68),
    @(248, 105, 67),
    @(249, 109, 43),
    @(250, 55, 37),
    @(251, 110, 91),
    @(252, 110, 84),
    @(253, 105, 45),
    @(254, 109, 70),
    @(255, 84, 51),
    @(256, "X", "X"),
    @(257, 120, 87),
    @(258, 107, 53),
    @(259, 136, 90),
    @(260, 65, 54),
    @(261, 107, 36),
    @(262, 60, 62),
    @(263, 30, 38),
    @(264, 129, 57),
    @(265, 111, 78),
    @(266, 110, 51),
    @(267, 30, 25),
    @(268, 85, 60),
    @(269, 155, 38),
    @(270, 61, 65),
    @(271, 92, 33),
    @(306, 45, 77),
    @(307, 42, 45),
    @(308, 11, 12),
    @(309, 45, 63),
    @(310, 87, 91),
    @(311, 93, 72),
    @(312, 75, 95),
    @(313, 71, 83),
    @(314, 72, 121),
    @(315, 85, 72),
    @(316, 52, 130),
    @(317, 67, 81),
    @(318, 47, 52)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 5).Value = $entry[1]
    $ws.Cells.Item($r, 6).Value = $entry[2]
}

# Move the view / selection to where the editing left off, and resize the
# window the way it was left in the saved workbook.
$win = $excel.ActiveWindow
$win.Left = 2860
$win.Top = 460
$win.Width = 24040
$win.Height = 16400
$win.ScrollRow = 247
$win.ScrollColumn = 1

$ws.Range("E271").Select()
